$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data, sorted descending by value, with Swedish & Uzbek removed
$data = @(
    @("Chinese", 20.84977367059123),
    @("English", 20.27228095637039),
    @("Spanish", 5.957929690022712),
    @("Arabic", 4.211816651393644),
    @("German", 3.819137072908771),
    @("Japanese", 3.450127318707894),
    @("Malay-Indonesian", 3.397911758346615),
    @("Russian", 2.808827039788575),
    @("Portuguese", 2.66637406592578),
    @("French", 2.390554103134028),
    @("Turkish", 2.054710512046594),
    @("Italian", 1.724640481899097),
    @("Korean", 1.640707820618253),
    @("Dutch", 1.126822919500077),
    @("Polish", 0.9801130171841317),
    @("Bengali", 0.948885184370871),
    @("Persian", 0.9407856270707152),
    @("Urdu", 0.9382682637495608),
    @("Vietnamese", 0.9236380623006191),
    @("Thai", 0.9182344117233137)
)

# Delete the two rows beyond the new data (currently rows go to 23; new data ends at row 21)
$ws.Range("A22:B23").Delete() | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
